# Update Correspond Handoff/Handback DateTime values for row 2 of the
# "zh-cn" and "de-de" worksheets to reflect a re-generated handback report.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 08:53:29"
$wsZhCn.Range("H2").Value = "2016-03-20 08:53:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 08:53:32"
$wsDeDe.Range("H2").Value = "2016-03-20 08:53:55"
